$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 112
$ws.Range("H112").Value = 1426.7115
$ws.Range("I112").Value = 383.66666
$ws.Range("J112").Value = 1562.7609
$ws.Range("K112").Value = 1150.99998
$ws.Range("L112").Value = 4688.2827
$ws.Range("M112").Value = -42.99998000000005
$ws.Range("N112").Value = -6904.2827

# Row 129
$ws.Range("H129").Value = 1097.5814
$ws.Range("J129").Value = 1050.5333
$ws.Range("L129").Value = 3151.5999
$ws.Range("N129").Value = -13151.5999

# Row 138
$ws.Range("H138").Value = 2888.9387
$ws.Range("J138").Value = 3552.3928
$ws.Range("L138").Value = 10657.1784
$ws.Range("N138").Value = -20937.1784

# ---------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 9956.9
$ws.Range("I32").Value = 9532.544
$ws.Range("K32").Value = 9532.544
$ws.Range("M32").Value = -9245.544

# Row 102
$ws.Range("H102").Value = 49235.453
$ws.Range("I102").Value = 16287.777
$ws.Range("J102").Value = 197500
$ws.Range("K102").Value = 16287.777
$ws.Range("L102").Value = 197500
$ws.Range("M102").Value = -14665.777
$ws.Range("N102").Value = -200744

# ---------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 105
$ws.Range("H105").Value = 2973.762
$ws.Range("I105").Value = 3215.3333
$ws.Range("J105").Value = 2369.8333
$ws.Range("K105").Value = 3215.3333
$ws.Range("L105").Value = 2369.8333
$ws.Range("M105").Value = -1468.3333
$ws.Range("N105").Value = -5863.8333

# Row 132
$ws.Range("H132").Value = 327867.7
$ws.Range("J132").Value = 1170710.6
$ws.Range("L132").Value = 3512131.8
$ws.Range("N132").Value = -3517191.8

# ---------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 68
$ws.Range("H68").Value = 1511.72
$ws.Range("J68").Value = 1643.2667
$ws.Range("L68").Value = 4929.800099999999
$ws.Range("N68").Value = -6551.800099999999

# Row 71
$ws.Range("H71").Value = 1511.72
$ws.Range("J71").Value = 1643.2667
$ws.Range("L71").Value = 14789.4003
$ws.Range("N71").Value = -22901.4003

# Rows 120-134 and 136-141: clear the computed-profit columns (H:N),
# leaving the leve metadata columns (A:G) intact. Row 135 is untouched.
$ws.Range("H120:N134").ClearContents()
$ws.Range("H136:N141").ClearContents()

# ---------------------------------------------------------------
# LTW sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 5880.8
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 5880.8
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 5880.8
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -6152.8

# ---------------------------------------------------------------
# WVR sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 41
$ws.Range("H41").Value = 21035.166
$ws.Range("J41").Value = 21035.166
$ws.Range("L41").Value = 21035.166
$ws.Range("N41").Value = -21815.166

# Row 45
$ws.Range("H45").Value = 11438.75
$ws.Range("I45").Value = 5564.5
$ws.Range("J45").Value = 17313
$ws.Range("K45").Value = 5564.5
$ws.Range("L45").Value = 17313
$ws.Range("M45").Value = -5073.5
$ws.Range("N45").Value = -18295
